# FueraPlazo.xlsx - "15022021 -modificacion a reporte solicitado desde
# gestion personas. -modificacion a excel para reporte"
#
# Changes applied:
#   1. Rename "Hoja1" -> "Cometido Fuera de Plazo" (also updates the
#      worksheet-scoped _FilterDatabase defined name automatically).
#   2. Add a new "Fecha Inicio" header in column E, matching the existing
#      header formatting/style used by A1 (bold, left/top aligned).
#   3. Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab. Excel automatically rewrites any
#    sheet-qualified references (e.g. the hidden _xlnm._FilterDatabase
#    defined name) to use the new name.
$ws.Name = "Cometido Fuera de Plazo"

# 2) New column E header: "Fecha Inicio". Copy A1's formatting first so the
#    new header cell reuses the same cell style as the other bold header
#    cells, then (re)write the value.
$ws.Range("E1").Value = "Fecha Inicio"
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "Fecha Inicio"

# 3) Freeze panes below row 1 (top row stays visible). Selecting the cell
#    that will become the top-left cell of the scrollable pane before
#    freezing matches Excel's own behaviour for a simple "freeze top row".
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
